# This workbook (Bmp8a-Tgfbr2 NATMI ligand-receptor output) was regenerated
# from an updated TPM expression matrix. Columns E:T (rows 2-17, one row per
# sending-cluster/target-cluster pair) are refreshed with the newly computed
# detection rates, expression values, specificities and edge weights; columns
# A:D (cluster/gene labels) and K:L are left as-is.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 16,16
# Row 2: ECs -> ECs
$arr[0,0] = 3.0
$arr[0,1] = 1.0
$arr[0,2] = 0.5725316666666668
$arr[0,3] = 1.717595
$arr[0,4] = 0.3864899584549088
$arr[0,5] = 0.3864899584549088
$arr[0,6] = 3.0
$arr[0,7] = 1.0
$arr[0,8] = 27.85106533333333
$arr[0,9] = 83.553196
$arr[0,10] = 0.1861900221007236
$arr[0,11] = 0.1861900221007236
$arr[0,12] = 15.94561685373556
$arr[0,13] = 143.51055168362
$arr[0,14] = 0.07196057390642722
$arr[0,15] = 0.07196057390642722
# Row 3: ECs -> FAPs
$arr[1,0] = 3.0
$arr[1,1] = 1.0
$arr[1,2] = 0.5725316666666668
$arr[1,3] = 1.717595
$arr[1,4] = 0.3864899584549088
$arr[1,5] = 0.3864899584549088
$arr[1,6] = 3.0
$arr[1,7] = 1.0
$arr[1,8] = 70.710031
$arr[1,9] = 212.130093
$arr[1,10] = 0.4727109026912454
$arr[1,11] = 0.4727109026912454
$arr[1,12] = 40.48373189848167
$arr[1,13] = 364.353587086335
$arr[1,14] = 0.1826980171423219
$arr[1,15] = 0.1826980171423219
# Row 4: ECs -> MuSCs
$arr[2,0] = 3.0
$arr[2,1] = 1.0
$arr[2,2] = 0.5725316666666668
$arr[2,3] = 1.717595
$arr[2,4] = 0.3864899584549088
$arr[2,5] = 0.3864899584549088
$arr[2,6] = 3.0
$arr[2,7] = 1.0
$arr[2,8] = 15.018964
$arr[2,9] = 45.056892
$arr[2,10] = 0.1004048213460311
$arr[2,11] = 0.1004048213460311
$arr[2,12] = 8.598832490526668
$arr[2,13] = 77.38949241474
$arr[2,14] = 0.03880545523070008
$arr[2,15] = 0.03880545523070008
# Row 5: ECs -> Resolving-Mac
$arr[3,0] = 3.0
$arr[3,1] = 1.0
$arr[3,2] = 0.5725316666666668
$arr[3,3] = 1.717595
$arr[3,4] = 0.3864899584549088
$arr[3,5] = 0.3864899584549088
$arr[3,6] = 3.0
$arr[3,7] = 1.0
$arr[3,8] = 36.00403133333333
$arr[3,9] = 108.012094
$arr[3,10] = 0.2406942538619999
$arr[3,11] = 0.2406942538619999
$arr[3,12] = 20.61344806599222
$arr[3,13] = 185.52103259393
$arr[3,14] = 0.09302591217545962
$arr[3,15] = 0.09302591217545962
# Row 6: FAPs -> ECs
$arr[4,0] = 2.0
$arr[4,1] = 0.6666666666666666
$arr[4,2] = 0.3064436666666667
$arr[4,3] = 0.919331
$arr[4,4] = 0.2068661122070742
$arr[4,5] = 0.2068661122070743
$arr[4,6] = 3.0
$arr[4,7] = 1.0
$arr[4,8] = 27.85106533333333
$arr[4,9] = 83.553196
$arr[4,10] = 0.1861900221007236
$arr[4,11] = 0.1861900221007236
$arr[4,12] = 8.534782581319556
$arr[4,13] = 76.813043231876
$arr[4,14] = 0.03851640600372592
$arr[4,15] = 0.03851640600372593
# Row 7: FAPs -> FAPs
$arr[5,0] = 2.0
$arr[5,1] = 0.6666666666666666
$arr[5,2] = 0.3064436666666667
$arr[5,3] = 0.919331
$arr[5,4] = 0.2068661122070742
$arr[5,5] = 0.2068661122070743
$arr[5,6] = 3.0
$arr[5,7] = 1.0
$arr[5,8] = 70.710031
$arr[5,9] = 212.130093
$arr[5,10] = 0.4727109026912454
$arr[5,11] = 0.4727109026912454
$arr[5,12] = 21.66864116975367
$arr[5,13] = 195.017770527783
$arr[5,14] = 0.09778786663763453
$arr[5,15] = 0.09778786663763453
# Row 8: FAPs -> MuSCs
$arr[6,0] = 2.0
$arr[6,1] = 0.6666666666666666
$arr[6,2] = 0.3064436666666667
$arr[6,3] = 0.919331
$arr[6,4] = 0.2068661122070742
$arr[6,5] = 0.2068661122070743
$arr[6,6] = 3.0
$arr[6,7] = 1.0
$arr[6,8] = 15.018964
$arr[6,9] = 45.056892
$arr[6,10] = 0.1004048213460311
$arr[6,11] = 0.1004048213460311
$arr[6,12] = 4.602466397694666
$arr[6,13] = 41.422197579252
$arr[6,14] = 0.0207703550386993
$arr[6,15] = 0.0207703550386993
# Row 9: FAPs -> Resolving-Mac
$arr[7,0] = 2.0
$arr[7,1] = 0.6666666666666666
$arr[7,2] = 0.3064436666666667
$arr[7,3] = 0.919331
$arr[7,4] = 0.2068661122070742
$arr[7,5] = 0.2068661122070743
$arr[7,6] = 3.0
$arr[7,7] = 1.0
$arr[7,8] = 36.00403133333333
$arr[7,9] = 108.012094
$arr[7,10] = 0.2406942538619999
$arr[7,11] = 0.2406942538619999
$arr[7,12] = 11.03320737656822
$arr[7,13] = 99.298866389114
$arr[7,14] = 0.04979148452701448
$arr[7,15] = 0.04979148452701449
# Row 10: MuSCs -> ECs
$arr[8,0] = 3.0
$arr[8,1] = 1.0
$arr[8,2] = 0.5451493333333333
$arr[8,3] = 1.635448
$arr[8,4] = 0.3680053968340403
$arr[8,5] = 0.3680053968340404
$arr[8,6] = 3.0
$arr[8,7] = 1.0
$arr[8,8] = 27.85106533333333
$arr[8,9] = 83.553196
$arr[8,10] = 0.1861900221007236
$arr[8,11] = 0.1861900221007236
$arr[8,12] = 15.18298969908978
$arr[8,13] = 136.646907291808
$arr[8,14] = 0.06851893296971553
$arr[8,15] = 0.06851893296971553
# Row 11: MuSCs -> FAPs
$arr[9,0] = 3.0
$arr[9,1] = 1.0
$arr[9,2] = 0.5451493333333333
$arr[9,3] = 1.635448
$arr[9,4] = 0.3680053968340403
$arr[9,5] = 0.3680053968340404
$arr[9,6] = 3.0
$arr[9,7] = 1.0
$arr[9,8] = 70.710031
$arr[9,9] = 212.130093
$arr[9,10] = 0.4727109026912454
$arr[9,11] = 0.4727109026912454
$arr[9,12] = 38.54752625962933
$arr[9,13] = 346.9277363366639
$arr[9,14] = 0.1739601633326692
$arr[9,15] = 0.1739601633326692
# Row 12: MuSCs -> MuSCs
$arr[10,0] = 3.0
$arr[10,1] = 1.0
$arr[10,2] = 0.5451493333333333
$arr[10,3] = 1.635448
$arr[10,4] = 0.3680053968340403
$arr[10,5] = 0.3680053968340404
$arr[10,6] = 3.0
$arr[10,7] = 1.0
$arr[10,8] = 15.018964
$arr[10,9] = 45.056892
$arr[10,10] = 0.1004048213460311
$arr[10,11] = 0.1004048213460311
$arr[10,12] = 8.187578211957332
$arr[10,13] = 73.68820390761599
$arr[10,14] = 0.03694951612349708
$arr[10,15] = 0.03694951612349709
# Row 13: MuSCs -> Resolving-Mac
$arr[11,0] = 3.0
$arr[11,1] = 1.0
$arr[11,2] = 0.5451493333333333
$arr[11,3] = 1.635448
$arr[11,4] = 0.3680053968340403
$arr[11,5] = 0.3680053968340404
$arr[11,6] = 3.0
$arr[11,7] = 1.0
$arr[11,8] = 36.00403133333333
$arr[11,9] = 108.012094
$arr[11,10] = 0.2406942538619999
$arr[11,11] = 0.2406942538619999
$arr[11,12] = 19.62757367867911
$arr[11,13] = 176.648163108112
$arr[11,14] = 0.08857678440815853
$arr[11,15] = 0.08857678440815854
# Row 14: Resolving-Mac -> ECs
$arr[12,0] = 2.0
$arr[12,1] = 0.6666666666666666
$arr[12,2] = 0.05723766666666667
$arr[12,3] = 0.171713
$arr[12,4] = 0.03863853250397663
$arr[12,5] = 0.03863853250397663
$arr[12,6] = 3.0
$arr[12,7] = 1.0
$arr[12,8] = 27.85106533333333
$arr[12,9] = 83.553196
$arr[12,10] = 0.1861900221007236
$arr[12,11] = 0.1861900221007236
$arr[12,12] = 1.594129993860889
$arr[12,13] = 14.347169944748
$arr[12,14] = 0.007194109220854936
$arr[12,15] = 0.007194109220854937
# Row 15: Resolving-Mac -> FAPs
$arr[13,0] = 2.0
$arr[13,1] = 0.6666666666666666
$arr[13,2] = 0.05723766666666667
$arr[13,3] = 0.171713
$arr[13,4] = 0.03863853250397663
$arr[13,5] = 0.03863853250397663
$arr[13,6] = 3.0
$arr[13,7] = 1.0
$arr[13,8] = 70.710031
$arr[13,9] = 212.130093
$arr[13,10] = 0.4727109026912454
$arr[13,11] = 0.4727109026912454
$arr[13,12] = 4.047277184367666
$arr[13,13] = 36.425494659309
$arr[13,14] = 0.01826485557861982
$arr[13,15] = 0.01826485557861982
# Row 16: Resolving-Mac -> MuSCs
$arr[14,0] = 2.0
$arr[14,1] = 0.6666666666666666
$arr[14,2] = 0.05723766666666667
$arr[14,3] = 0.171713
$arr[14,4] = 0.03863853250397663
$arr[14,5] = 0.03863853250397663
$arr[14,6] = 3.0
$arr[14,7] = 1.0
$arr[14,8] = 15.018964
$arr[14,9] = 45.056892
$arr[14,10] = 0.1004048213460311
$arr[14,11] = 0.1004048213460311
$arr[14,12] = 0.8596504551106666
$arr[14,13] = 7.736854095996
$arr[14,14] = 0.003879494953134587
$arr[14,15] = 0.003879494953134588
# Row 17: Resolving-Mac -> Resolving-Mac
$arr[15,0] = 2.0
$arr[15,1] = 0.6666666666666666
$arr[15,2] = 0.05723766666666667
$arr[15,3] = 0.171713
$arr[15,4] = 0.03863853250397663
$arr[15,5] = 0.03863853250397663
$arr[15,6] = 3.0
$arr[15,7] = 1.0
$arr[15,8] = 36.00403133333333
$arr[15,9] = 108.012094
$arr[15,10] = 0.2406942538619999
$arr[15,11] = 0.2406942538619999
$arr[15,12] = 2.060786744113555
$arr[15,13] = 18.547080697022
$arr[15,14] = 0.009300072751367286
$arr[15,15] = 0.009300072751367287

$ws.Range("E2:T17").Value = $arr

Write-Host "Updated E2:T17 with recalculated TPM-based NATMI values."
